# Use timezone from preferences for excel reports
# Replace the joda-time based date formatting expressions used in the
# route export template with calls to dateTool.format(...), which takes
# the locale/timezone from preferences.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'
$ws.Range("B9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", position.fixTime, locale, timezone)}'

$ws.Range("B2").Select()
